$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: remove the empty B46 inline-string cell entirely (dropped in the diff)
$ws.Range("B46").ClearContents()

# Row 47
$ws.Range('A47').Value = "'120"
$ws.Range('B47').Value = "'AbnedAK"
$ws.Range('C47').Value = 0
$ws.Range('D47').Value = 0
$ws.Range('E47').Value = 0
$ws.Range('F47').Value = 0
$ws.Range('G47').Value = 0
$ws.Range('H47').Value = 0
$ws.Range('I47').Value = 1
$ws.Range('J47').Value = 0
$ws.Range('K47').Value = 0
$ws.Range('L47').Value = 0
$ws.Range('M47').Value = 0
$ws.Range('N47').Value = 1
$ws.Range('O47').Value = 0
$ws.Range('P47').Value = 0
$ws.Range('Q47').Value = 0
$ws.Range('R47').Value = 1
$ws.Range('S47').Value = 0
$ws.Range('T47').Value = 0
$ws.Range('U47').Value = 1
$ws.Range('V47').Value = 0
$ws.Range('W47').Value = 0
$ws.Range('X47').Value = 0
$ws.Range('Y47').Value = 0
$ws.Range('Z47').Value = 0
$ws.Range('AA47').Value = 0
$ws.Range('AB47').Value = 0
$ws.Range('AC47').Value = 0
$ws.Range('AD47').Value = 0
$ws.Range('AE47').Value = 0
$ws.Range('AF47').Value = 0

# Row 48
$ws.Range('A48').Value = "'11110"
$ws.Range('B48').Value = "'MatkoudAN"
$ws.Range('C48').Value = 1
$ws.Range('D48').Value = 0
$ws.Range('E48').Value = 0
$ws.Range('F48').Value = 0
$ws.Range('G48').Value = 0
$ws.Range('H48').Value = 1
$ws.Range('I48').Value = 0
$ws.Range('J48').Value = 0
$ws.Range('K48').Value = 0
$ws.Range('L48').Value = 0
$ws.Range('M48').Value = 0
$ws.Range('N48').Value = 0
$ws.Range('O48').Value = 0
$ws.Range('P48').Value = 1
$ws.Range('Q48').Value = 0
$ws.Range('R48').Value = 0
$ws.Range('S48').Value = 0
$ws.Range('T48').Value = 0
$ws.Range('U48').Value = 0
$ws.Range('V48').Value = 0
$ws.Range('W48').Value = 0
$ws.Range('X48').Value = 0
$ws.Range('Y48').Value = 0
$ws.Range('Z48').Value = 0
$ws.Range('AA48').Value = 0
$ws.Range('AB48').Value = 0
$ws.Range('AC48').Value = 0
$ws.Range('AD48').Value = 0
$ws.Range('AE48').Value = 1
$ws.Range('AF48').Value = 0

# Row 49
$ws.Range('A49').Value = "'110"
$ws.Range('C49').Value = 1
$ws.Range('D49').Value = 0
$ws.Range('E49').Value = 1
$ws.Range('F49').Value = 1
$ws.Range('G49').Value = 0
$ws.Range('H49').Value = 0
$ws.Range('I49').Value = 1
$ws.Range('J49').Value = 0
$ws.Range('K49').Value = 0
$ws.Range('L49').Value = 0
$ws.Range('M49').Value = 0
$ws.Range('N49').Value = 0
$ws.Range('O49').Value = 0
$ws.Range('P49').Value = 0
$ws.Range('Q49').Value = 0
$ws.Range('R49').Value = 1
$ws.Range('S49').Value = 0
$ws.Range('T49').Value = 0
$ws.Range('U49').Value = 0
$ws.Range('V49').Value = 0
$ws.Range('W49').Value = 0
$ws.Range('X49').Value = 0
$ws.Range('Y49').Value = 0
$ws.Range('Z49').Value = 1
$ws.Range('AA49').Value = 0
$ws.Range('AB49').Value = 0
$ws.Range('AC49').Value = 0
$ws.Range('AD49').Value = 0
$ws.Range('AE49').Value = 0
$ws.Range('AF49').Value = 0

# Row 50
$ws.Range('A50').Value = "'2012"
$ws.Range('C50').Value = 0
$ws.Range('D50').Value = 1
$ws.Range('E50').Value = 1
$ws.Range('F50').Value = 0
$ws.Range('G50').Value = 1
$ws.Range('H50').Value = 0
$ws.Range('I50').Value = 0
$ws.Range('J50').Value = 0
$ws.Range('K50').Value = 0
$ws.Range('L50').Value = 1
$ws.Range('M50').Value = 0
$ws.Range('N50').Value = 1
$ws.Range('O50').Value = 0
$ws.Range('P50').Value = 0
$ws.Range('Q50').Value = 0
$ws.Range('R50').Value = 1
$ws.Range('S50').Value = 0
$ws.Range('T50').Value = 0
$ws.Range('U50').Value = 0
$ws.Range('V50').Value = 1
$ws.Range('W50').Value = 0
$ws.Range('X50').Value = 0
$ws.Range('Y50').Value = 0
$ws.Range('Z50').Value = 0
$ws.Range('AA50').Value = 0
$ws.Range('AB50').Value = 0
$ws.Range('AC50').Value = 0
$ws.Range('AD50').Value = 0
$ws.Range('AE50').Value = 0
$ws.Range('AF50').Value = 0

# Row 51
$ws.Range('A51').Value = "'2"
$ws.Range('C51').Value = 1
$ws.Range('D51').Value = 1
$ws.Range('E51').Value = 1
$ws.Range('F51').Value = 1
$ws.Range('G51').Value = 0
$ws.Range('H51').Value = 0
$ws.Range('I51').Value = 0
$ws.Range('J51').Value = 0
$ws.Range('K51').Value = 0
$ws.Range('L51').Value = 0
$ws.Range('M51').Value = 0
$ws.Range('N51').Value = 0
$ws.Range('O51').Value = 0
$ws.Range('P51').Value = 0
$ws.Range('Q51').Value = 0
$ws.Range('R51').Value = 0
$ws.Range('S51').Value = 0
$ws.Range('T51').Value = 0
$ws.Range('U51').Value = 0
$ws.Range('V51').Value = 0
$ws.Range('W51').Value = 0
$ws.Range('X51').Value = 0
$ws.Range('Y51').Value = 0
$ws.Range('Z51').Value = 0
$ws.Range('AA51').Value = 0
$ws.Range('AB51').Value = 0
$ws.Range('AC51').Value = 1
$ws.Range('AD51').Value = 0
$ws.Range('AE51').Value = 0
$ws.Range('AF51').Value = 0

# Row 52
$ws.Range('A52').Value = "'19320"
$ws.Range('C52').Value = 0
$ws.Range('D52').Value = 0
$ws.Range('E52').Value = 0
$ws.Range('F52').Value = 0
$ws.Range('G52').Value = 0
$ws.Range('H52').Value = 0
$ws.Range('I52').Value = 1
$ws.Range('J52').Value = 0
$ws.Range('K52').Value = 0
$ws.Range('L52').Value = 0
$ws.Range('M52').Value = 0
$ws.Range('N52').Value = 0
$ws.Range('O52').Value = 0
$ws.Range('P52').Value = 0
$ws.Range('Q52').Value = 0
$ws.Range('R52').Value = 0
$ws.Range('S52').Value = 0
$ws.Range('T52').Value = 0
$ws.Range('U52').Value = 0
$ws.Range('V52').Value = 0
$ws.Range('W52').Value = 0
$ws.Range('X52').Value = 0
$ws.Range('Y52').Value = 0
$ws.Range('Z52').Value = 0
$ws.Range('AA52').Value = 0
$ws.Range('AB52').Value = 0
$ws.Range('AC52').Value = 0
$ws.Range('AD52').Value = 0
$ws.Range('AE52').Value = 0

# Row 53
$ws.Range('A53').Value = "'1019"
$ws.Range('C53').Value = 0
$ws.Range('D53').Value = 0
$ws.Range('E53').Value = 0
$ws.Range('F53').Value = 0
$ws.Range('G53').Value = 1
$ws.Range('H53').Value = 0
$ws.Range('I53').Value = 0
$ws.Range('J53').Value = 0
$ws.Range('K53').Value = 0
$ws.Range('L53').Value = 0
$ws.Range('M53').Value = 0
$ws.Range('N53').Value = 0
$ws.Range('O53').Value = 1
$ws.Range('P53').Value = 0
$ws.Range('Q53').Value = 1
$ws.Range('R53').Value = 0
$ws.Range('S53').Value = 1
$ws.Range('T53').Value = 0
$ws.Range('U53').Value = 0
$ws.Range('V53').Value = 0
$ws.Range('W53').Value = 0
$ws.Range('X53').Value = 0
$ws.Range('Y53').Value = 0
$ws.Range('Z53').Value = 0
$ws.Range('AA53').Value = 0
$ws.Range('AB53').Value = 1
$ws.Range('AC53').Value = 0
$ws.Range('AD53').Value = 1
$ws.Range('AE53').Value = 0

# Row 54
$ws.Range('A54').Value = "'1092113"
$ws.Range('C54').Value = 0
$ws.Range('D54').Value = 1
$ws.Range('E54').Value = 1
$ws.Range('F54').Value = 0
$ws.Range('G54').Value = 1
$ws.Range('H54').Value = 0
$ws.Range('I54').Value = 0
$ws.Range('J54').Value = 0
$ws.Range('K54').Value = 0
$ws.Range('L54').Value = 1
$ws.Range('M54').Value = 0
$ws.Range('N54').Value = 1
$ws.Range('O54').Value = 0
$ws.Range('P54').Value = 0
$ws.Range('Q54').Value = 0
$ws.Range('R54').Value = 1
$ws.Range('S54').Value = 0
$ws.Range('T54').Value = 0
$ws.Range('U54').Value = 0
$ws.Range('V54').Value = 1
$ws.Range('W54').Value = 0
$ws.Range('X54').Value = 0
$ws.Range('Y54').Value = 0
$ws.Range('Z54').Value = 0
$ws.Range('AA54').Value = 0
$ws.Range('AB54').Value = 0
$ws.Range('AC54').Value = 0
$ws.Range('AD54').Value = 0
$ws.Range('AE54').Value = 0
$ws.Range('AF54').Value = 0

# Row 55
$ws.Range('A55').Value = "'1092113"
$ws.Range('C55').Value = 0
$ws.Range('D55').Value = 1
$ws.Range('E55').Value = 1
$ws.Range('F55').Value = 0
$ws.Range('G55').Value = 1
$ws.Range('H55').Value = 0
$ws.Range('I55').Value = 0
$ws.Range('J55').Value = 0
$ws.Range('K55').Value = 0
$ws.Range('L55').Value = 1
$ws.Range('M55').Value = 0
$ws.Range('N55').Value = 1
$ws.Range('O55').Value = 0
$ws.Range('P55').Value = 0
$ws.Range('Q55').Value = 0
$ws.Range('R55').Value = 1
$ws.Range('S55').Value = 0
$ws.Range('T55').Value = 0
$ws.Range('U55').Value = 0
$ws.Range('V55').Value = 1
$ws.Range('W55').Value = 0
$ws.Range('X55').Value = 0
$ws.Range('Y55').Value = 0
$ws.Range('Z55').Value = 0
$ws.Range('AA55').Value = 0
$ws.Range('AB55').Value = 0
$ws.Range('AC55').Value = 0
$ws.Range('AD55').Value = 0
$ws.Range('AE55').Value = 0
$ws.Range('AF55').Value = 0

# Row 56
$ws.Range('A56').Value = "'0"
$ws.Range('C56').Value = 1
$ws.Range('D56').Value = 0
$ws.Range('E56').Value = 1
$ws.Range('F56').Value = 0
$ws.Range('G56').Value = 1
$ws.Range('H56').Value = 0
$ws.Range('I56').Value = 0
$ws.Range('J56').Value = 0
$ws.Range('K56').Value = 0
$ws.Range('L56').Value = 0
$ws.Range('M56').Value = 0
$ws.Range('N56').Value = 0
$ws.Range('O56').Value = 1
$ws.Range('P56').Value = 0
$ws.Range('Q56').Value = 0

# Row 57
$ws.Range('A57').Value = "'0"
$ws.Range('C57').Value = 1
$ws.Range('D57').Value = 0
$ws.Range('E57').Value = 1
$ws.Range('F57').Value = 0
$ws.Range('G57').Value = 1
$ws.Range('H57').Value = 0
$ws.Range('I57').Value = 0
$ws.Range('J57').Value = 0
$ws.Range('K57').Value = 0
$ws.Range('L57').Value = 0
$ws.Range('M57').Value = 0
$ws.Range('N57').Value = 0
$ws.Range('O57').Value = 1
$ws.Range('P57').Value = 0
$ws.Range('Q57').Value = 0

# Row 58
$ws.Range('A58').Value = "'19320"
$ws.Range('C58').Value = 1
$ws.Range('D58').Value = 0
$ws.Range('E58').Value = 1
$ws.Range('F58').Value = 1
$ws.Range('G58').Value = 0
$ws.Range('H58').Value = 0
$ws.Range('I58').Value = 1
$ws.Range('J58').Value = 0
$ws.Range('K58').Value = 0
$ws.Range('L58').Value = 0
$ws.Range('M58').Value = 0
$ws.Range('N58').Value = 0
$ws.Range('O58').Value = 0
$ws.Range('P58').Value = 0
$ws.Range('Q58').Value = 0
$ws.Range('R58').Value = 1
$ws.Range('S58').Value = 0
$ws.Range('T58').Value = 0
$ws.Range('U58').Value = 0
$ws.Range('V58').Value = 0
$ws.Range('W58').Value = 0
$ws.Range('X58').Value = 0
$ws.Range('Y58').Value = 0
$ws.Range('Z58').Value = 1
$ws.Range('AA58').Value = 0
$ws.Range('AB58').Value = 0
$ws.Range('AC58').Value = 0
$ws.Range('AD58').Value = 0
$ws.Range('AE58').Value = 0
$ws.Range('AF58').Value = 0

# Row 59
$ws.Range('A59').Value = "'000"
$ws.Range('C59').Value = 0
$ws.Range('D59').Value = 0
$ws.Range('E59').Value = 0
$ws.Range('F59').Value = 0
$ws.Range('G59').Value = 0
$ws.Range('H59').Value = 0
$ws.Range('I59').Value = 0
$ws.Range('J59').Value = 0
$ws.Range('K59').Value = 0
$ws.Range('L59').Value = 0
$ws.Range('M59').Value = 0
$ws.Range('N59').Value = 0
$ws.Range('O59').Value = 0
$ws.Range('P59').Value = 0

# Row 60
$ws.Range('A60').Value = "'000"
$ws.Range('C60').Value = 0
$ws.Range('D60').Value = 0
$ws.Range('E60').Value = 0
$ws.Range('F60').Value = 0
$ws.Range('G60').Value = 0
$ws.Range('H60').Value = 0
$ws.Range('I60').Value = 0
$ws.Range('J60').Value = 0
$ws.Range('K60').Value = 0
$ws.Range('L60').Value = 0
$ws.Range('M60').Value = 0
$ws.Range('N60').Value = 0
$ws.Range('O60').Value = 0
$ws.Range('P60').Value = 0

# Row 61
$ws.Range('A61').Value = "'000"
$ws.Range('C61').Value = 0
$ws.Range('D61').Value = 0
$ws.Range('E61').Value = 0
$ws.Range('F61').Value = 0
$ws.Range('G61').Value = 0
$ws.Range('H61').Value = 0
$ws.Range('I61').Value = 0
$ws.Range('J61').Value = 0
$ws.Range('K61').Value = 0
$ws.Range('L61').Value = 0
$ws.Range('M61').Value = 0
$ws.Range('N61').Value = 0
$ws.Range('O61').Value = 0
$ws.Range('P61').Value = 0

# Row 62
$ws.Range('A62').Value = "'09099"
$ws.Range('C62').Value = 0
$ws.Range('D62').Value = 0
$ws.Range('E62').Value = 0
$ws.Range('F62').Value = 0
$ws.Range('G62').Value = 0
$ws.Range('H62').Value = 0
$ws.Range('I62').Value = 0
$ws.Range('J62').Value = 0
$ws.Range('K62').Value = 0
$ws.Range('L62').Value = 0
$ws.Range('M62').Value = 0
$ws.Range('N62').Value = 0
$ws.Range('O62').Value = 0
$ws.Range('P62').Value = 1

# Row 63
$ws.Range('A63').Value = "'09099"
$ws.Range('C63').Value = 0
$ws.Range('D63').Value = 0
$ws.Range('E63').Value = 0
$ws.Range('F63').Value = 0
$ws.Range('G63').Value = 0
$ws.Range('H63').Value = 0
$ws.Range('I63').Value = 0
$ws.Range('J63').Value = 0
$ws.Range('K63').Value = 0
$ws.Range('L63').Value = 0
$ws.Range('M63').Value = 0
$ws.Range('N63').Value = 0
$ws.Range('O63').Value = 0
$ws.Range('P63').Value = 1

# Row 64
$ws.Range('A64').Value = "'2303"
$ws.Range('C64').Value = 0
$ws.Range('D64').Value = 1
$ws.Range('E64').Value = 0
$ws.Range('F64').Value = 1
$ws.Range('G64').Value = 0
$ws.Range('H64').Value = 0
$ws.Range('I64').Value = 0
$ws.Range('J64').Value = 0
$ws.Range('K64').Value = 0
$ws.Range('L64').Value = 0
$ws.Range('M64').Value = 0
$ws.Range('N64').Value = 0
$ws.Range('O64').Value = 0
$ws.Range('P64').Value = 0

# Row 65
$ws.Range('A65').Value = "'10301"
$ws.Range('B65').Value = "'MohamedAhied"
$ws.Range('C65').Value = 1
$ws.Range('D65').Value = 0
$ws.Range('E65').Value = 0
$ws.Range('F65').Value = 1
$ws.Range('G65').Value = 0
$ws.Range('H65').Value = 0
$ws.Range('I65').Value = 0
$ws.Range('J65').Value = 0
$ws.Range('K65').Value = 0
$ws.Range('L65').Value = 0
$ws.Range('M65').Value = 0
$ws.Range('N65').Value = 0
$ws.Range('O65').Value = 0
$ws.Range('P65').Value = 0
$ws.Range('Q65').Value = 1

# Row 66
$ws.Range('A66').Value = "'100"
$ws.Range('B66').Value = "'"
$ws.Range('C66').Value = 0
$ws.Range('D66').Value = 0
$ws.Range('E66').Value = 0
$ws.Range('F66').Value = 0
$ws.Range('G66').Value = 0
$ws.Range('H66').Value = 0
$ws.Range('I66').Value = 0
$ws.Range('J66').Value = 0
$ws.Range('K66').Value = 0
$ws.Range('L66').Value = 0
$ws.Range('M66').Value = 0
$ws.Range('N66').Value = 0
$ws.Range('O66').Value = 0
$ws.Range('P66').Value = 0
